# Adding sex selection to weight for age plot example.
$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet ---
# Insert a new row at 46: an intro note before the age/weight questions.
$survey.Rows.Item(46).Insert()
$survey.Range("A46").Value = "note"
$survey.Range("D46").Value = "The following data will be used to generate a weight for age plot."

# Existing age question, now shifted to row 47: tweak its label & add a hint.
$survey.Range("D47").Value = "Enter age (in years):"
$survey.Range("E47").Value = "Must be less than 20."

# Existing weight question, now shifted to row 48: tweak its label.
$survey.Range("D48").Value = "Enter weight (in lbs):"

# Insert a new row at 49: the new sex question, before "end screen".
$survey.Rows.Item(49).Insert()
$survey.Range("A49").Value = "select_one sexes"
$survey.Range("C49").Value = "sex"
$survey.Range("D49").Value = "Enter sex:"

# --- choices sheet ---
# Add the new "sexes" choice list (male / female) after the existing choices.
$choices.Range("A23").Value = "sexes"
$choices.Range("B23").Value = "male"
$choices.Range("C23").Value = "male"

$choices.Range("A24").Value = "sexes"
$choices.Range("B24").Value = "female"
$choices.Range("C24").Value = "female"
